$wb = $excel.ActiveWorkbook

# ---- Settings sheet: new Asset Type rows 28-32 ----
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A28").Value = 'VS'
$wsSettings.Range("B28").Value = 'Verification Services'
$wsSettings.Range("A29").Value = 'SEC'
$wsSettings.Range("B29").Value = 'Self-Employment Certification'
$wsSettings.Range("A30").Value = 1040
$wsSettings.Range("B30").Value = 1040
$wsSettings.Range("A31").Value = 'ScheduleC'
$wsSettings.Range("B31").Value = 'Schedule C'
$wsSettings.Range("A32").Value = 'ProfitLoss'
$wsSettings.Range("B32").Value = 'Profit/Loss Statement'

# ---- Assets sheet: new selector row 18 ----
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("A18").Value = 'GetTextUploadPageNum_Selector'
$wsAssets.Range("B18").Value = 'GetTextUploadPageNum_Selector'

# ---- Findings sheet: insert PS_SumOfEarningsMatch row + append WN/VS/SEC rule blocks ----
$wsFindings = $wb.Worksheets.Item("Findings")
$wsFindings.Rows.Item(102).Insert()
$wsFindings.Range("A102").Value = 'PS_SumOfEarningsMatch'
$wsFindings.Range("B102").Value = 'Gross Pay listed on the Paystub does not match the sum of the earnings on the Paystub.'

# WN_* additional rules (rows 138-142)
$wsFindings.Range("A138").Value = 'WN_TotalMatchSum'
$wsFindings.Range("B138").Value = 'Total Amount listed on the Work Number does not match the sum of the base, overtime, commission, bonus, and other columns.'
$wsFindings.Range("A139").Value = 'WN_PayPeriodEndingExistsInEDC'
$wsFindings.Range("B139").Value = 'Two most recent Period Endings listed on the Work Number and Exact Day Calculator do not match.'
$wsFindings.Range("A140").Value = 'WN_PayPeriodStartingExistsInEDC'
$wsFindings.Range("B140").Value = 'Two most recent Period Starts listed on the Work Number and Exact Day Calculator do not match.'
$wsFindings.Range("A141").Value = 'WN_RecentGrossPayMatch'
$wsFindings.Range("B141").Value = 'Two most recent Gross Earnings listed on the Work Number and Exact Day Calculator do not match.'
$wsFindings.Range("A142").Value = 'WN_HighestCalValue'
$wsFindings.Range("B142").Value = 'Highest Calculated Income on the Exact Day Calculator does not match Verified Income on the ICW.'

# VS (Verification Services) rules (rows 144-152, row 143 stays blank separator)
$wsFindings.Range("A144").Value = 'VS_Datecheck'
$wsFindings.Range("B144").Value = 'Date listed on Verification Services is over 120 days in the past.'
$wsFindings.Range("A145").Value = 'VS_EmployerCheck'
$wsFindings.Range("B145").Value = 'Employer Name listed on Verification Services, Exact Day Calculator, and ICW do not match.'
$wsFindings.Range("A146").Value = 'VS_EmployeeCheck'
$wsFindings.Range("B146").Value = 'Employee Name listed on Verification Services, Exact Day Calculator, and ICW do not match.'
$wsFindings.Range("A147").Value = 'VS_PayFrequency'
$wsFindings.Range("B147").Value = 'Pay Frequency listed on Verification Services, Exact Day Calculator, and ICW do not match.'
$wsFindings.Range("A148").Value = 'VS_PerYearEarningsCheck'
$wsFindings.Range("B148").Value = 'The application is for The Safford property but the Year to Date Income is not complete in the Exact Day Calculator.'
$wsFindings.Range("A149").Value = 'VS_PayPeriodStartingExistsInEDC'
$wsFindings.Range("B149").Value = 'Two most recent Period Starts listed on the Verification Services and Exact Day Calculator do not match.'
$wsFindings.Range("A150").Value = 'VS_PayPeriodEndingExistsInEDC'
$wsFindings.Range("B150").Value = 'Two most recent Period Endings listed on the Verification Services and Exact Day Calculator do not match.'
$wsFindings.Range("A151").Value = 'VS_RecentGrossPayMatch'
$wsFindings.Range("B151").Value = 'Two most recent Gross Earnings listed on the Verification Services and Exact Day Calculator do not match.'
$wsFindings.Range("A152").Value = 'VS_HighestCalValue'
$wsFindings.Range("B152").Value = 'Highest Calculated Income on the Exact Day Calculator does not match Verified Income on the ICW.'

# SEC (Self-Employment Certification) rules (rows 154-167, row 153 stays blank separator)
$wsFindings.Range("A154").Value = 'SEC_AllFieldsDocumented'
$wsFindings.Range("B154").Value = 'All fields were not documented appropriately.'
$wsFindings.Range("A155").Value = 'SEC_CheckDate'
$wsFindings.Range("B155").Value = 'Date listed on the Self-Employment Certification is over 120 days in the past.'
$wsFindings.Range("A156").Value = 'SEC_EmployeeCheck'
$wsFindings.Range("B156").Value = 'Name listed on the Self-Employment Certification and ICW do not match.'
$wsFindings.Range("A157").Value = 'SEC_AmountCheck'
$wsFindings.Range("B157").Value = 'Amount listed on the Self-Employment Certification and ICW do not match.'
$wsFindings.Range("A158").Value = 'SEC_IDoCheck'
$wsFindings.Range("B158").Value = 'The "I do" box was checked on the Self-Employment Certification. Sending for manual verification. '
$wsFindings.Range("A159").Value = 'SEC_ConfirmSupportingDocuments'
$wsFindings.Range("B159").Value = 'Could not confirm appropriate supporting documents (such as tax documents or profit/loss statement) were included in the application.'
$wsFindings.Range("A160").Value = 'SEC_NameListCheck'
$wsFindings.Range("B160").Value = 'Name listed on the Self-Employment Certification, 1040, and Schedule C do not match.'
$wsFindings.Range("A161").Value = 'SEC_AmountListCheck'
$wsFindings.Range("B161").Value = 'Amount listed on the Self-Employment Certification, 1040, and Schedule C do not line up as expected. Needs additional review. '
$wsFindings.Range("A162").Value = 'SEC_1040NotSigned'
$wsFindings.Range("B162").Value = 'The 1040 tax form was not signed.'
$wsFindings.Range("A163").Value = 'SEC_1040NotDated'
$wsFindings.Range("B163").Value = 'The 1040 tax form was not dated appropriately.'
$wsFindings.Range("A164").Value = 'SEC_NameOnProfitLoss'
$wsFindings.Range("B164").Value = 'Name listed on the Self-Employment Certification and Profit/Lost Statement do not match.'
$wsFindings.Range("A165").Value = 'SEC_NetIncomeonProfitLoss'
$wsFindings.Range("B165").Value = 'Net Income displayed on the Profit/Loss Statement doe not match the calculated Net Income.'
$wsFindings.Range("A166").Value = 'SEC_AmountOnProfitLoss'
$wsFindings.Range("B166").Value = 'Amount listed on the Self-Employment Certification and Profit/Loss Statement do not line up as expected. Needs additional review. '
$wsFindings.Range("A167").Value = 'SEC_ExpensesAreGeneric'
$wsFindings.Range("B167").Value = 'Expenses listed are too generic, need more specific expense names.'

# ---- Selection / view state (best-effort cosmetic match) ----
$wsSettings.Range("A33").Select()
$wsAssets.Range("B24").Select()
$wsFindings.Range("A166").Select()
$wsFindings.Activate()
